$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.Style = $origStyle
}

Set-TextValue 'D2' '90.407.63'
Set-TextValue 'E2' '  -0.39%  '
Set-TextValue 'D3' '3.092.67'
Set-TextValue 'E3' '  -1.94%  '
Set-TextValue 'E4' '  +0.19%  '
Set-TextValue 'D5' '234.33'
Set-TextValue 'E5' '  +8.60%  '
Set-TextValue 'D6' '622.95'
Set-TextValue 'E6' '  -0.53%  '
Set-TextValue 'E7' '  -4.21%  '
Set-TextValue 'D8' '0.365'
Set-TextValue 'E8' '  -1.14%  '
Set-TextValue 'D10' '3.090.20'
Set-TextValue 'E10' '  -1.83%  '
Set-TextValue 'D11' '0.728'
Set-TextValue 'E11' '  -3.69%  '
Set-TextValue 'D12' '0.197'
Set-TextValue 'E12' '  -1.75%  '
Set-TextValue 'B13' 'ShibaInu'
Set-TextValue 'C13' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D13' '0.0000251'
Set-TextValue 'E13' '  +2.76%  '
Set-TextValue 'B14' 'Avalanche'
Set-TextValue 'C14' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D14' '36.18'
Set-TextValue 'E14' '  +4.14%  '
Set-TextValue 'D15' '5.47'
Set-TextValue 'E15' '  -2.91%  '
Set-TextValue 'D16' '90.022.70'
Set-TextValue 'E16' '  -0.35%  '
Set-TextValue 'D17' '3.660.18'
Set-TextValue 'E17' '  -1.84%  '
Set-TextValue 'D18' '3.077.77'
Set-TextValue 'E18' '  -3.29%  '
Set-TextValue 'E19' '  +5.12%  '
Set-TextValue 'D20' '0.0000217'
Set-TextValue 'E20' '  +3.32%  '
Set-TextValue 'D21' '14.01'
Set-TextValue 'E21' '  -1.92%  '
Set-TextValue 'B22' 'BitcoinCash'
Set-TextValue 'C22' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D22' '436.55'
Set-TextValue 'E22' '  -5.35%  '
Set-TextValue 'B23' 'Polkadot'
Set-TextValue 'C23' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D23' '5.57'
Set-TextValue 'E23' '  +5.55%  '
Set-TextValue 'D24' '8.90'
Set-TextValue 'E24' '  -1.49%  '
Set-TextValue 'D25' '5.94'
Set-TextValue 'E25' '  +1.34%  '
Set-TextValue 'D26' '7.54'
Set-TextValue 'E26' '  -1.90%  '
Set-TextValue 'D27' '88.79'
Set-TextValue 'E27' '  -0.06%  '
Set-TextValue 'D28' '12.10'
Set-TextValue 'E28' '  -0.08%  '
Set-TextValue 'D29' '3.253.51'
Set-TextValue 'E29' '  -1.80%  '
Set-TextValue 'D30' '1.00'
Set-TextValue 'E30' '  +0.00%  '
Set-TextValue 'D31' '9.37'
Set-TextValue 'E31' '  +1.79%  '
Set-TextValue 'E32' '  -1.88%  '
Set-TextValue 'E33' '  -2.85%  '
Set-TextValue 'D34' '0.198'
Set-TextValue 'E34' '  +9.71%  '
Set-TextValue 'E35' '  +9.56%  '
Set-TextValue 'D36' '25.83'
Set-TextValue 'E36' '  -4.85%  '
Set-TextValue 'D37' '3.85'
Set-TextValue 'E37' '  +6.50%  '
Set-TextValue 'D38' '503.41'
Set-TextValue 'E38' '  -2.82%  '
Set-TextValue 'D39' '7.16'
Set-TextValue 'E39' '  +4.22%  '
Set-TextValue 'E40' '  -0.30%  '
Set-TextValue 'D41' '1.29'
Set-TextValue 'E41' '  +0.58%  '
Set-TextValue 'D42' '0.0887'
Set-TextValue 'E42' '  +1.66%  '
Set-TextValue 'D43' '22.16'
Set-TextValue 'E43' '  -0.22%  '
Set-TextValue 'D44' '0.405'
Set-TextValue 'E44' '  -1.68%  '
Set-TextValue 'D46' '3.44'
Set-TextValue 'E46' '  +55.48%  '
Set-TextValue 'E47' '  -1.87%  '
Set-TextValue 'D48' '0.692'
Set-TextValue 'E48' '  +3.99%  '
Set-TextValue 'D49' '152.57'
Set-TextValue 'E49' '  +3.22%  '
Set-TextValue 'D50' '44.88'
Set-TextValue 'E50' '  -3.16%  '
Set-TextValue 'E51' '  +0.11%  '
